$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.158893346786499
$ws.Range("B1").Value = 2.666262865066528
$ws.Range("C1").Value = 3.563746452331543
$ws.Range("D1").Value = 6.073071956634521
$ws.Range("E1").Value = 2.055503606796265
